$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New registrant rows appended below the existing data (rows 14-15),
# matching the "DOB / PhoneNumber / SendText" columns that are stored as
# literal text elsewhere in this sheet (e.g. row 8-13), and the
# "CurrentDate / FutureDate" columns that use the custom date/time
# number format already present in the workbook.

$dateFmt = "yyyy-mm-dd h:mm:ss"

# ---- Row 14: Jason / J ----
$ws.Range("A14").Value = "Jason"
$ws.Range("C14").Value = "J"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "100193"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1111111111"

$ws.Range("F14").Value = "jason@gmail.com"
$ws.Range("G14").Value = "one Washington"

$ws.Range("H14").NumberFormat = $dateFmt
$ws.Range("H14").Value = 44610.00347222222

$ws.Range("I14").NumberFormat = $dateFmt
$ws.Range("I14").Value = 44631

$ws.Range("J14").Value = 1

$ws.Range("K14").NumberFormat = "@"
$ws.Range("K14").Value = "+11111111111"

# ---- Row 15: Pranav / Konduru ----
$ws.Range("A15").Value = "Pranav"
$ws.Range("C15").Value = "Konduru"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "10100000"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "9255777863"

$ws.Range("F15").Value = "pranav.konduru@sjsu.edu"
$ws.Range("G15").Value = "One Washington Square"

$ws.Range("H15").NumberFormat = $dateFmt
$ws.Range("H15").Value = 44610.00347222222

$ws.Range("I15").NumberFormat = $dateFmt
$ws.Range("I15").Value = 44631

$ws.Range("J15").Value = 1

$ws.Range("K15").NumberFormat = "@"
$ws.Range("K15").Value = "+19255777863"
